$wb = $excel.ActiveWorkbook

foreach ($name in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 432
    $ws.Range("F3").Value = 7
}
